# Update the two-digit multiplication problems in the single table of the
# document. The table has 5 columns; only every 5th row (1, 5, 10, 15, 20)
# actually holds problem text, the rest are blank "work space" rows.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$rows = @(1, 5, 10, 15, 20)

$values = @(
    @("94×39=", "22×25=", "30×16=", "52×94=", "82×96="),
    @("93×18=", "12×55=", "63×66=", "95×86=", "98×16="),
    @("36×37=", "84×93=", "31×95=", "70×73=", "15×71="),
    @("48×26=", "56×47=", "17×30=", "11×67=", "54×30="),
    @("71×20=", "81×21=", "58×54=", "50×57=", "83×65=")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$i][$c - 1]
    }
}
